$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Mordialloc/Cheltenham old row replaced with a standalone exposure-time row
$ws.Range("A2").Value = "30/12/20 10:45am- 12:15pm"
$ws.Range("B2:D2").ClearContents()
$ws.Range("E2").Value = "new"

# Row 3: Bentleigh
$ws.Range("A3").Value = "Bentleigh"
$ws.Range("B3").Value = "Il Centro Deli  5/284/292 Centre Rd, Bentleigh VIC 3204"
$ws.Range("C3").Value = "22/12/20 12:00pm-12:30pm"
$ws.Range("D3").Value = "Case shopped in store"
$ws.Range("E3").Value = "new"

# Row 4: Brighton
$ws.Range("A4").Value = "Brighton"
$ws.Range("B4").Value = "Brighton Beach  Brighton, VIC 3186"
$ws.Range("C4").Value = "29/12/20 12:00pm-3:00pm"
$ws.Range("D4").Value = "Case attended beach"
$ws.Range("E4").Value = "new"

# Row 5: Cheltenham - Angus and Cootes Jeweller
$ws.Range("A5").Value = "Cheltenham"
$ws.Range("B5").Value = "Angus and Cootes Jeweller  Southland Shopping Centre, 2096/1239 Nepean Hwy, Cheltenham VIC 3192"
$ws.Range("C5").Value = "28/12/2020 2:30pm-2:50pm"
$ws.Range("D5").Value = "Case shopped in store"
$ws.Range("E5").Value = "new"

# Row 6: Cheltenham - Honey Birdette Southland
$ws.Range("A6").Value = "Cheltenham"
$ws.Range("B6").Value = "Honey Birdette Southland  Shop 2209/1239, Southland Shopping Centre, Cheltenham VIC 3192"
$ws.Range("C6").Value = "22/12/2020 3:50pm-4:05pm"
$ws.Range("D6").Value = "Case shopped in store"
$ws.Range("E6").Value = "new"

# Row 7: Cheltenham - Mecca Southland
$ws.Range("A7").Value = "Cheltenham"
$ws.Range("B7").Value = "Mecca Southland  Shop 2011/2013, Southland Shopping Centre, Cheltenham VIC 3192"
$ws.Range("C7").Value = "22/12/2020 3:30pm-3:50pm"
$ws.Range("D7").Value = "Case shopped in store"
$ws.Range("E7").Value = "new"

# Row 8: Mentone (new)
$ws.Range("A8").Value = "Mentone"
$ws.Range("B8").Value = "Woolworths Mentone  105-111 Balcombe Road, Mentone VIC 3194"
$ws.Range("C8").Value = "23/12/20 2:45pm-3:05pm"
$ws.Range("D8").Value = "Case shopped in store"
$ws.Range("E8").Value = "new"

# Row 9: Moorabbin (new)
$ws.Range("A9").Value = "Moorabbin"
$ws.Range("B9").Value = "COSTCO Moorabbin  8 Chifley Drive, Moorabbin Airport VIC 3194"
$ws.Range("C9").Value = "30/12/20 4:00m- 5:50pm"
$ws.Range("D9").Value = "Case shopped in store"
$ws.Range("E9").Value = "new"

# Column width adjustments to match new content (closest achievable values in this runtime)
$ws.Columns.Item(1).ColumnWidth = 23.17
$ws.Columns.Item(2).ColumnWidth = 82.17
$ws.Columns.Item(4).ColumnWidth = 17.17

# Selection changed to full-column selection A:E
[void]$ws.Range("A1:E1048576").Select()
